$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("InsurancePremium")

# --- Columns G, M, T (Vehicle_Date of Manufacture / Insurant_birthdate / Product_startdate) ---
# switch from numeric date serials to plain text (dd/mm/yyyy) values, formatted as Text ("@").
$ws.Columns("G").NumberFormat = "@"
$ws.Columns("M").NumberFormat = "@"
$ws.Columns("T").NumberFormat = "@"

# Row 2 - Honda MotorCycle
$ws.Range("G2").Value = "04/09/2019"
$ws.Range("M2").Value = "12/12/1989"
$ws.Range("T2").Value = "15/06/2023"
$ws.Range("U2").Value = 3000000

# Row 3 - Volvo Scooter
$ws.Range("G3").Value = "04/09/2001"
$ws.Range("M3").Value = "12/12/1977"
$ws.Range("T3").Value = "15/06/2023"
$ws.Range("U3").Value = 3000000

# Row 4 - Audi Moped
$ws.Range("G4").Value = "04/09/1989"
$ws.Range("M4").Value = "12/12/1956"
$ws.Range("T4").Value = "15/06/2023"
$ws.Range("U4").Value = 3000000

# --- Leftover formatted-but-empty cells below the table (rows 7-11), same date format as col G/M/T originally used ---
$ws.Range("S7").NumberFormat = "dd\/mm\/yyyy"

$ws.Range("L8").NumberFormat = "dd\/mm\/yyyy"
$ws.Range("S8").NumberFormat = "dd\/mm\/yyyy"

$ws.Range("F9").NumberFormat = "dd\/mm\/yyyy"
$ws.Range("L9").NumberFormat = "dd\/mm\/yyyy"
$ws.Range("S9").NumberFormat = "dd\/mm\/yyyy"

$ws.Range("F10").NumberFormat = "dd\/mm\/yyyy"
$ws.Range("L10").NumberFormat = "dd\/mm\/yyyy"

$ws.Range("F11").NumberFormat = "dd\/mm\/yyyy"

# Selection settles back on A1 (no explicit multi-cell selection remains in saved file)
$ws.Range("A1").Select()
